# Applies the "Remoção de linhas inutilizadas" edit:
#   - Adds a new column F ("Dividas:") that holds the last "<Categoria> Situação:"
#     section name found in column E, wrapped as [['Categoria']].
#   - Strips the trailing "TOTAL ORIGEM: ..." line from column E and removes the
#     "TOTAL:" marker at the end of each sub-total line (keeping the trailing
#     space before the newline), leaving a trailing newline at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (copy format from E1).
$ws.Range("F1").Value = "Dividas:"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

$lastRow = 144

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $text = $eCell.Value2

    if ($text -eq $null) { continue }

    # Derive the "Dividas:" category from the last "<...> Situação:" header line.
    $headerMatches = [regex]::Matches($text, "(?m)^(.*) Situação:$")
    $lastCategory = $headerMatches[$headerMatches.Count - 1].Groups[1].Value
    $fValue = "[['" + $lastCategory + "']]"

    # Strip the trailing "TOTAL ORIGEM: ..." line.
    $marker = "`nTOTAL ORIGEM:"
    $cut = $text.LastIndexOf($marker)
    if ($cut -ge 0) {
        $newText = $text.Substring(0, $cut)
    } else {
        $newText = $text
    }

    # Remove the "TOTAL:" marker at the end of each running-total line, keeping
    # the space that preceded it, then terminate with a trailing newline.
    $newText = $newText.Replace(" TOTAL:", " ") + "`n"

    $eCell.Value = $newText
    $ws.Cells.Item($r, 6).Value = $fValue
}

Write-Host "Applied Dividas column + trimmed TOTAL ORIGEM lines for rows 2..$lastRow"
